# "changes done In Admin Create Page Test"
# Update the Start/Expected/Actual Completion date values on row 2 of the
# Project_Details sheet, and move the active selection to M8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start Date (L2), Expected Date (M2), Actual Completion Date (N2) -
# stored as date serial numbers (cells already carry the dd/mmm/yyyy style).
$ws.Range("L2").Value = 37201
$ws.Range("M2").Value = 48357
$ws.Range("N2").Value = 48226

# Move/save the current selection on the sheet.
$ws.Range("M8").Select()
